$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column D ("Edate" - event end date) between Sdate (C) and Cimage (old D, now E).
# Excel's Insert() shifts existing columns D..I to E..J and carries formatting
# (column width / bestFit) from the column to the left, which matches the target layout.
$ws.Columns.Item(4).Insert()

# Give the new Edate column a width close to the adjacent Sdate column (cosmetic).
$ws.Columns.Item(4).ColumnWidth = 15.5

# Header for the new column
$ws.Range("D1").Value = "Edate"

# Populate Edate values.
# Football, Cricket and Arm Wrestling now have an explicit end date.
$ws.Range("D2").Value = "2023-11-08T00:00:00"
$ws.Range("D3").Value = "2023-11-08T00:00:00"
$ws.Range("D4").Value = "2023-11-08T00:00:00"

# The remaining events (no dates set yet) get a numeric placeholder of 0,
# matching their Sdate (column C) placeholder value.
$ws.Range("D5").Value = 0
$ws.Range("D6").Value = 0
$ws.Range("D7").Value = 0
$ws.Range("D8").Value = 0

# Football and Cricket (rows 2 and 3) move from UPCOMING to ONGOING.
# Columns after the insert: F=UPCOMING, G=ONGOING, H=COMPLETED.
$ws.Range("F2").Value = 0
$ws.Range("G2").Value = 1
$ws.Range("H2").Value = 0

$ws.Range("F3").Value = 0
$ws.Range("G3").Value = 1
$ws.Range("H3").Value = 0
